# "Add files via upload" -- fills in the previously-blank "计划内容"
# (plan content) column for the 第七周周四 (week 7, Thursday) task block,
# rows 33-38 on Sheet1, with the "设计界面原型[...]" (design UI prototype)
# entries for each team member, and updates the sheet's active
# selection/scroll position to reflect where the user ended up editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B33").Value = "设计界面原型[手机端用户注册]"

# This row's cell picked up the explicit "宋体" font style (cellXfs s="5")
# rather than the theme-linked default used by the rest of the column.
$ws.Range("B34").Value = "设计界面原型[手机端用户登陆]"
$ws.Range("B34").Font.Name = "宋体"

$ws.Range("B35").Value = "设计界面原型[PC端用户注册]"
$ws.Range("B36").Value = "设计界面原型[PC端用户登陆]"
$ws.Range("B37").Value = "设计界面原型[登录管理系统]"
$ws.Range("B38").Value = "设计界面原型[查看用户信息]"

# Reflect the author's final cursor/scroll position (A16 at the top of the
# viewport, B38 selected) from the saved sheetView.
$ws.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

# Minor cosmetic tab-area/scrollbar split the author also nudged while
# editing (bookViews/workbookView tabRatio 560 -> 550).
$excel.ActiveWindow.TabRatio = 0.55
